$d = $word.ActiveDocument

# Actions, keyed by the ORIGINAL (pre-edit) 1-based paragraph index.
# Because we always process from the highest index down to the lowest,
# deleting/restyling a paragraph never disturbs the index of any
# paragraph we have not yet processed (paragraph indices in this
# engine are re-resolved against the live collection on each access,
# so we must never hold a stale Paragraph reference across an earlier
# mutation - we always re-fetch by index immediately before use).
#
# "delete"            -> remove the whole paragraph
# "invisibleseparator" -> change style to invisibleseparator and text to " "
# @("style", X)        -> change paragraph style to X, leave text alone

$actions = @{
    1  = @("delete")
    2  = @("delete")
    3  = @("delete")
    4  = @("style", "SectionNumber")
    5  = @("style", "SectionTitle")
    6  = @("delete")
    7  = @("style", "SectionSubtitle")
    8  = @("delete")
    9  = @("style", "SectionHeadnote")
    10 = @("invisibleseparator")
    19 = @("delete")
    20 = @("delete")
    21 = @("delete")
    24 = @("delete")
    38 = @("delete")
    39 = @("delete")
    40 = @("invisibleseparator")
    43 = @("delete")
    45 = @("delete")
    46 = @("delete")
    47 = @("invisibleseparator")
    50 = @("delete")
    52 = @("invisibleseparator")
    53 = @("delete")
}

$maxIndex = 53
for ($i = $maxIndex; $i -ge 1; $i--) {
    if (-not $actions.ContainsKey($i)) {
        continue
    }
    $action = $actions[$i]
    $p = $d.Paragraphs($i)

    if ($action[0] -eq "delete") {
        $p.Range.Delete()
    } elseif ($action[0] -eq "invisibleseparator") {
        $p.Range.Text = " "
        $p.Style = "invisibleseparator"
    } elseif ($action[0] -eq "style") {
        $p.Style = $action[1]
    }
}

# --- Remove now-unused style definitions ---
$stylesToDelete = @("Chapter Spacer", "Head End", "Head Field Separator", "Head Separator", "Node End", "Node Start")
foreach ($styleName in $stylesToDelete) {
    $s = $d.Styles($styleName)
    $s.Delete()
}

Write-Output "done"
